# Archive.xlsx — "Update functions and Data Model" (#50)
#
# 1. Noémi Villars's surname is updated to "Noémi Villars-Amberg" in the
#    existing Authorship cell (H2).
# 2. A new "Authorship Resource" column (I) is added, carrying the same
#    author string as column H.
# 3. The stale "general alignment" style that used to sit on most of the
#    data row is cleared back to the workbook's default ("Normal") style —
#    only the Directory cell (D2) keeps the bordered/left-aligned style.
# 4. The saved selection is moved to C15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the author's surname -------------------------------------------------
$ws.Range("H2").Value = "Daniela Subotic, Noémi Villars-Amberg"

# --- 2. Add the "Authorship Resource" column --------------------------------------
# Clone the header's look (bordered, left-aligned) onto the new header cell,
# then fill in the header text and the row-2 value (mirrors column H).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Authorship Resource"
$ws.Range("I2").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("I1").ColumnWidth = 44.6667

# --- 3. Drop the unused "general" style from row 2 (except Directory) -------------
foreach ($addr in @("A2", "B2", "C2", "E2", "F2", "G2", "H2")) {
    $ws.Range($addr).Style = "Normal"
}

# --- 4. Restore the saved selection -------------------------------------------
[void]$ws.Range("C15").Select()
